$wb = $excel.ActiveWorkbook

# --- Update header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 83
$wsForecast.Range("C2").Value = -37.598013276473
$wsForecast.Range("D2").Value = 198.4826571779394
$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("B3").Value = 83
$wsForecast.Range("C3").Value = -38.75267517789888
$wsForecast.Range("D3").Value = 198.504523264877
$wsForecast.Range("A4").Value = 44948.99999999999
$wsForecast.Range("B4").Value = 83
$wsForecast.Range("C4").Value = -35.83055927495901
$wsForecast.Range("D4").Value = 206.2603017451686
$wsForecast.Range("A5").Value = 44962.99999999999
$wsForecast.Range("B5").Value = 82
$wsForecast.Range("C5").Value = -38.76598761553932
$wsForecast.Range("D5").Value = 202.9694688290953
$wsForecast.Range("A6").Value = 44969.99999999999
$wsForecast.Range("B6").Value = 82
$wsForecast.Range("C6").Value = -42.42006607119447
$wsForecast.Range("D6").Value = 206.6267709308632
$wsForecast.Range("A7").Value = 44976.99999999999
$wsForecast.Range("B7").Value = 82
$wsForecast.Range("C7").Value = -41.3796109025008
$wsForecast.Range("D7").Value = 200.8252620154346
$wsForecast.Range("A8").Value = 44983.99999999999
$wsForecast.Range("B8").Value = 81
$wsForecast.Range("C8").Value = -42.83911814353748
$wsForecast.Range("D8").Value = 199.1649416291972
$wsForecast.Range("A9").Value = 44990.99999999999
$wsForecast.Range("B9").Value = 81
$wsForecast.Range("C9").Value = -39.74107032072772
$wsForecast.Range("D9").Value = 199.8013098069874
$wsForecast.Range("A10").Value = 44997.99999999999
$wsForecast.Range("B10").Value = 81
$wsForecast.Range("C10").Value = -40.50189431064364
$wsForecast.Range("D10").Value = 198.7612958361314
$wsForecast.Range("A11").Value = 45011.99999999999
$wsForecast.Range("B11").Value = 80
$wsForecast.Range("C11").Value = -37.98760203058615
$wsForecast.Range("D11").Value = 200.4717826190727
$wsForecast.Range("A12").Value = 45018.99999999999
$wsForecast.Range("B12").Value = 80
$wsForecast.Range("C12").Value = -37.77483497188051
$wsForecast.Range("D12").Value = 202.1577945605025
$wsForecast.Range("A13").Value = 45032.99999999999
$wsForecast.Range("B13").Value = 80
$wsForecast.Range("C13").Value = -39.13667808907554
$wsForecast.Range("D13").Value = 191.9632377118342
$wsForecast.Range("A14").Value = 45039.99999999999
$wsForecast.Range("B14").Value = 79
$wsForecast.Range("C14").Value = -35.14322738247283
$wsForecast.Range("D14").Value = 191.0562602393443
$wsForecast.Range("A15").Value = 45046.99999999999
$wsForecast.Range("B15").Value = 79
$wsForecast.Range("C15").Value = -45.67996763414252
$wsForecast.Range("D15").Value = 197.5594622018089
$wsForecast.Range("A16").Value = 45053.99999999999
$wsForecast.Range("B16").Value = 79
$wsForecast.Range("C16").Value = -43.75119094279339
$wsForecast.Range("D16").Value = 198.5636464773883
$wsForecast.Range("A17").Value = 45088.99999999999
$wsForecast.Range("B17").Value = 78
$wsForecast.Range("C17").Value = -47.82259251591258
$wsForecast.Range("D17").Value = 197.2307791019536
$wsForecast.Range("A18").Value = 45095.99999999999
$wsForecast.Range("B18").Value = 77
$wsForecast.Range("C18").Value = -37.98467427588691
$wsForecast.Range("D18").Value = 193.8262814241448
$wsForecast.Range("A19").Value = 45102.99999999999
$wsForecast.Range("B19").Value = 77
$wsForecast.Range("C19").Value = -48.41462493629822
$wsForecast.Range("D19").Value = 201.0109622178131
$wsForecast.Range("A20").Value = 45109.99999999999
$wsForecast.Range("B20").Value = 77
$wsForecast.Range("C20").Value = -40.75228819313531
$wsForecast.Range("D20").Value = 200.799913677403
$wsForecast.Range("A21").Value = 45116.99999999999
$wsForecast.Range("B21").Value = 77
$wsForecast.Range("C21").Value = -40.97646986640805
$wsForecast.Range("D21").Value = 189.7120549061759
$wsForecast.Range("A22").Value = 45123.99999999999
$wsForecast.Range("B22").Value = 76
$wsForecast.Range("C22").Value = -47.48241442008103
$wsForecast.Range("D22").Value = 194.540376832205
$wsForecast.Range("A23").Value = 45130.99999999999
$wsForecast.Range("B23").Value = 76
$wsForecast.Range("C23").Value = -44.88983551128739
$wsForecast.Range("D23").Value = 195.3710776757362
$wsForecast.Range("A24").Value = 45137.99999999999
$wsForecast.Range("B24").Value = 76
$wsForecast.Range("C24").Value = -52.64334928121934
$wsForecast.Range("D24").Value = 190.0069892139585
$wsForecast.Range("A25").Value = 45144.99999999999
$wsForecast.Range("B25").Value = 76
$wsForecast.Range("C25").Value = -42.57838145192445
$wsForecast.Range("D25").Value = 191.6843797211634
$wsForecast.Range("A26").Value = 45151.99999999999
$wsForecast.Range("B26").Value = 75
$wsForecast.Range("C26").Value = -47.71728419617654
$wsForecast.Range("D26").Value = 189.4202013767274
$wsForecast.Range("A27").Value = 45158.99999999999
$wsForecast.Range("B27").Value = 75
$wsForecast.Range("C27").Value = -43.75094667086312
$wsForecast.Range("D27").Value = 195.1195966230337
$wsForecast.Range("A28").Value = 45165.99999999999
$wsForecast.Range("B28").Value = 75
$wsForecast.Range("C28").Value = -43.62812036521606
$wsForecast.Range("D28").Value = 194.9776998293016
$wsForecast.Range("A29").Value = 45347.99999999999
$wsForecast.Range("B29").Value = 69
$wsForecast.Range("C29").Value = -50.40021454631362
$wsForecast.Range("D29").Value = 192.8766639507254
$wsForecast.Range("A30").Value = 45354.99999999999
$wsForecast.Range("B30").Value = 68
$wsForecast.Range("C30").Value = -40.6150783265244
$wsForecast.Range("D30").Value = 184.3988885265454
$wsForecast.Range("A31").Value = 45361.99999999999
$wsForecast.Range("B31").Value = 68
$wsForecast.Range("C31").Value = -51.67399116865662
$wsForecast.Range("D31").Value = 192.9807068872273
$wsForecast.Range("A32").Value = 45368.99999999999
$wsForecast.Range("B32").Value = 68
$wsForecast.Range("C32").Value = -60.66393689933863
$wsForecast.Range("D32").Value = 187.6707835890621
$wsForecast.Range("A33").Value = 45382.99999999999
$wsForecast.Range("B33").Value = 67
$wsForecast.Range("C33").Value = -47.30048523146942
$wsForecast.Range("D33").Value = 189.8927486407972
$wsForecast.Range("A34").Value = 45389.99999999999
$wsForecast.Range("B34").Value = 67
$wsForecast.Range("C34").Value = -50.90871791372717
$wsForecast.Range("D34").Value = 183.1648117754334
$wsForecast.Range("A35").Value = 45396.99999999999
$wsForecast.Range("B35").Value = 67
$wsForecast.Range("C35").Value = -55.46768125725465
$wsForecast.Range("D35").Value = 185.4832584885398
$wsForecast.Range("A36").Value = 45480.99999999999
$wsForecast.Range("B36").Value = 64
$wsForecast.Range("C36").Value = -60.55433179007139
$wsForecast.Range("D36").Value = 179.9171411593573
$wsForecast.Range("A37").Value = 45494.99999999999
$wsForecast.Range("B37").Value = 63
$wsForecast.Range("C37").Value = -50.29541098754665
$wsForecast.Range("D37").Value = 182.3181606174111
$wsForecast.Range("A38").Value = 45501.99999999999
$wsForecast.Range("B38").Value = 63
$wsForecast.Range("C38").Value = -60.90557009568913
$wsForecast.Range("D38").Value = 187.4976004948539
$wsForecast.Range("A39").Value = 45515.99999999999
$wsForecast.Range("B39").Value = 63
$wsForecast.Range("C39").Value = -46.03549630304442
$wsForecast.Range("D39").Value = 185.7629322521187
$wsForecast.Range("A40").Value = 45529.99999999999
$wsForecast.Range("B40").Value = 62
$wsForecast.Range("C40").Value = -57.04779604066113
$wsForecast.Range("D40").Value = 186.5403309190393
$wsForecast.Range("A41").Value = 45536.99999999999
$wsForecast.Range("B41").Value = 62
$wsForecast.Range("C41").Value = -62.17562752813383
$wsForecast.Range("D41").Value = 177.9171161999453
$wsForecast.Range("A42").Value = 45543.99999999999
$wsForecast.Range("B42").Value = 62
$wsForecast.Range("C42").Value = -60.1082109340759
$wsForecast.Range("D42").Value = 173.9091244743662
$wsForecast.Range("A43").Value = 45550.99999999999
$wsForecast.Range("B43").Value = 61
$wsForecast.Range("C43").Value = -51.38017451579248
$wsForecast.Range("D43").Value = 190.3559271876168
$wsForecast.Range("A44").Value = 45557.99999999999
$wsForecast.Range("B44").Value = 61
$wsForecast.Range("C44").Value = -60.24703219651794
$wsForecast.Range("D44").Value = 179.5577127063544
$wsForecast.Range("A45").Value = 45571.99999999999
$wsForecast.Range("B45").Value = 61
$wsForecast.Range("C45").Value = -63.53953253384415
$wsForecast.Range("D45").Value = 186.0144593239524
$wsForecast.Range("A46").Value = 45578.99999999999
$wsForecast.Range("B46").Value = 60
$wsForecast.Range("C46").Value = -59.06793979824822
$wsForecast.Range("D46").Value = 165.2841113866724
$wsForecast.Range("A47").Value = 45585.99999999999
$wsForecast.Range("B47").Value = 60
$wsForecast.Range("C47").Value = -60.32438232145171
$wsForecast.Range("D47").Value = 179.5825228250941
$wsForecast.Range("A48").Value = 45592.99999999999
$wsForecast.Range("B48").Value = 60
$wsForecast.Range("C48").Value = -62.44777555760031
$wsForecast.Range("D48").Value = 174.4586470905649
$wsForecast.Range("A49").Value = 45599.99999999999
$wsForecast.Range("B49").Value = 60
$wsForecast.Range("C49").Value = -63.03369746841155
$wsForecast.Range("D49").Value = 187.8996505012414
$wsForecast.Range("A50").Value = 45606.99999999999
$wsForecast.Range("B50").Value = 59
$wsForecast.Range("C50").Value = -56.73335262522812
$wsForecast.Range("D50").Value = 184.2670947115664
$wsForecast.Range("A51").Value = 45613.99999999999
$wsForecast.Range("B51").Value = 59
$wsForecast.Range("C51").Value = -65.85498960308095
$wsForecast.Range("D51").Value = 183.842517136718
$wsForecast.Range("A52").Value = 45620.99999999999
$wsForecast.Range("B52").Value = 59
$wsForecast.Range("C52").Value = -60.30270940763586
$wsForecast.Range("D52").Value = 172.9197181185741
$wsForecast.Range("A53").Value = 45627.99999999999
$wsForecast.Range("B53").Value = 59
$wsForecast.Range("C53").Value = -63.20177062556795
$wsForecast.Range("D53").Value = 179.2493289956711
$wsForecast.Range("A54").Value = 45634.99999999999
$wsForecast.Range("B54").Value = 58
$wsForecast.Range("C54").Value = -63.79943341780032
$wsForecast.Range("D54").Value = 181.1528914589485
$wsForecast.Range("A55").Value = 45641.99999999999
$wsForecast.Range("B55").Value = 58
$wsForecast.Range("C55").Value = -54.11729946975956
$wsForecast.Range("D55").Value = 182.7271411164385
$wsForecast.Range("A56").Value = 45648.99999999999
$wsForecast.Range("B56").Value = 58
$wsForecast.Range("C56").Value = -61.75475085682893
$wsForecast.Range("D56").Value = 168.6204607399958
$wsForecast.Range("A57").Value = 45655.99999999999
$wsForecast.Range("B57").Value = 58
$wsForecast.Range("C57").Value = -63.65129624936704
$wsForecast.Range("D57").Value = 180.3747538753896
$wsForecast.Range("A58").Value = 45662.99999999999
$wsForecast.Range("B58").Value = 57
$wsForecast.Range("C58").Value = -62.22006106135753
$wsForecast.Range("D58").Value = 177.3664300753771
$wsForecast.Range("A59").Value = 45669.99999999999
$wsForecast.Range("B59").Value = 57
$wsForecast.Range("C59").Value = -62.15020258390854
$wsForecast.Range("D59").Value = 183.1357283488026
$wsForecast.Range("A60").Value = 45676.99999999999
$wsForecast.Range("B60").Value = 57
$wsForecast.Range("C60").Value = -66.00777825775923
$wsForecast.Range("D60").Value = 173.0235543356748
$wsForecast.Range("A61").Value = 45683.99999999999
$wsForecast.Range("B61").Value = 57
$wsForecast.Range("C61").Value = -61.11293905838006
$wsForecast.Range("D61").Value = 174.6304854884061

# --- Apply formatting to match source workbook conventions ---
# Header style (bold, bordered, centered) copied from an existing header cell
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Date-serial number-format style copied from an existing date column cell
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the originally active sheet/selection
$wsWeekly.Activate()
$null = $wsWeekly.Range("A1").Select()

Write-Output "PO Forecast sheet created with $($wsForecast.UsedRange.Rows.Count) rows"
